$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing table one column to the right (A:L -> B:M) to make
# room for a new leading "Scenario" column.
$ws.Columns("A:A").Insert()

# New "Java-Large" dataset header/divider row, copying the same
# (bold header / dark-fill divider) formatting used for the Java14 block.
$ws.Range("B6").Value = "Java-Large"
$ws.Range("C2:M2").Copy()
$ws.Range("C6:M6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Scenario labels for the existing Java14 block.
$ws.Range("A3").Value = "Scenario 1"
$ws.Range("A4").Value = "Scenario 3"
$ws.Range("A5").Value = "Scenario 3"

# New Java-Large data rows (Scenario 1: Original, Scenario 3: Krakatau/Procyon).
$ws.Range("A7").Value = "Scenario 1"
$ws.Range("B7").Value = "Original"
$ws.Range("C7").Value = "N/A"
$ws.Range("D7").Value = 87
$ws.Range("E7").Value = 344
$ws.Range("F7").Value = 8988
$ws.Range("G7").Value = 12
$ws.Range("H7").Value = 656
$ws.Range("I7").Value = 0.4148
$ws.Range("J7").Value = 0.9663
$ws.Range("K7").Value = 0.344
$ws.Range("L7").Value = 0.8301
$ws.Range("M7").Value = 0.5074

$ws.Range("A8").Value = "Scenario 3"
$ws.Range("B8").Value = "Krakatau"
$ws.Range("C8").Value = "N/A"
$ws.Range("D8").Value = 91
$ws.Range("E8").Value = 628
$ws.Range("F8").Value = 9000
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 372
$ws.Range("I8").Value = 0.4279
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = 0.628
$ws.Range("L8").Value = 0.8738
$ws.Range("M8").Value = 0.7715

$ws.Range("A9").Value = "Scenario 3"
$ws.Range("B9").Value = "Procyon"
$ws.Range("C9").Value = "N/A"
$ws.Range("D9").Value = 88
$ws.Range("E9").Value = 788
$ws.Range("F9").Value = 8976
$ws.Range("G9").Value = 24
$ws.Range("H9").Value = 212
$ws.Range("I9").Value = 0.4339
$ws.Range("J9").Value = 0.9704
$ws.Range("K9").Value = 0.788
$ws.Range("L9").Value = 0.9208
$ws.Range("M9").Value = 0.8698

$ws.Range("L9").Select()
